# Cotações atualizadas - 2025-12-11
# Append the new quote row (97) below the existing data, mirroring the
# date/number formatting of the previous row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 96's formatting down to row 97 first so the new date cell
# inherits the same date/time number format (style index 2) as the rest
# of column A.
$ws.Range("A96").Copy($ws.Range("A97"))

$ws.Range("A97").Value = 46002
$ws.Range("B97").Value = "22,0482"
$ws.Range("C97").Value = "15,9866"
$ws.Range("D97").Value = "15,5979"
$ws.Range("E97").Value = "15,5979"
